$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 64.2
$ws.Range("N2").Value = 85.8724807945396

$ws.Range("K3").Value = 54.8
$ws.Range("N3").Value = 85.8724807945396

$ws.Range("K4").Value = 53
$ws.Range("N4").Value = 85.8724807945396

$ws.Range("K5").Value = 52.8
$ws.Range("N5").Value = 85.8724807945396

$ws.Range("K6").Value = 49
$ws.Range("N6").Value = 85.8724807945396
